# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a freshly generated
# handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-27 13:00:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-27 13:00:47"
$zhcn.Range("K2").Value = "2016-08-27 13:01:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-27 13:01:27"
